# Applies the "plate num" sheet update:
#  - Fills in D:G (sex/environ/sire/dam) values for data rows 2-7
#  - Removes the now-unused reference rows 16-18 (p1=limnetic / hatfield / plate num)
#  - Updates the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sex/environ/sire/dam values for each data row
$values = @(
    @{ Row = 2; D = "U"; E = 1; F = 1; G = 1 },
    @{ Row = 3; D = "U"; E = 1; F = 2; G = 2 },
    @{ Row = 4; D = "U"; E = 1; F = 2; G = 1 },
    @{ Row = 5; D = "U"; E = 1; F = 3; G = 3 },
    @{ Row = 6; D = "U"; E = 1; F = 3; G = 1 },
    @{ Row = 7; D = "U"; E = 1; F = 3; G = 2 }
)

foreach ($item in $values) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

# Remove the old reference rows (16-18) that held "p1=limnetic", "hatfield", "plate num"
$ws.Range("A16:G18").Value = $null

# Update the selected cell to match the target state
$ws.Range("G13").Select()
